# Insert a new price-history snapshot column before the "nom" column (ER),
# pushing "nom" -> ES and "url_produit" -> ET, matching the commit:
# "Update LDLC prices history".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column ER (148) is currently "nom". Insert a new blank column there;
# Excel shifts ER->ES and ES->ET automatically (formats/shared content too).
$ws.Columns("ER").Insert()

# Determine the extent of data (last used row) on the sheet.
$usedRows = $ws.UsedRange.Rows.Count

# New header for the freshly inserted ER column: the latest timestamp snapshot.
$ws.Cells.Item(1, 148).Value = "2026-02-03 17:43:21"

# For every data row, copy the most recent snapshot price (now in column EQ,
# col 147) into the new ER column (148) -- but only where that latest
# snapshot actually holds a price (rows where the product had gone out of
# the tracked/priced state stay blank, same as column EQ).
for ($r = 2; $r -le $usedRows; $r++) {
    $latest = $ws.Cells.Item($r, 147).Value()
    if ($latest -ne $null -and $latest -ne "") {
        $ws.Cells.Item($r, 148).Value = $latest
    }
}
